$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows above the existing item rows (old rows 7,8 -> new rows 9,10)
$ws.Rows("7:8").Insert()

# Copy formatting (styles + merges) from the rows that used to be the data rows
# (now shifted to 9 and 10) onto the two freshly inserted blank rows.
$ws.Range("A9:Q9").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)
$ws.Range("A10:Q10").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)

# Match row heights of the template rows.
$ws.Rows(7).RowHeight = $ws.Rows(9).RowHeight
$ws.Rows(8).RowHeight = $ws.Rows(10).RowHeight

# --- Row 7: new item "ALPHINTERN 30 F.C.TABS" ---
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "ALPHINTERN 30 F.C.TABS"
$ws.Range("H7").Value = "0:2"
$ws.Range("L7").Value = "1"
$ws.Range("N7").Value = "87.00"
$ws.Range("P7").Value = "28.7100"
$ws.Range("Q7").Value = "0:1"

# --- Row 8: new item "EXTRAUMA DNA FORTE TOPICAL GEL 25 GM" ---
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "EXTRAUMA DNA FORTE TOPICAL GEL 25 GM"
$ws.Range("H8").Value = "1:0"
$ws.Range("L8").Value = "1"
$ws.Range("N8").Value = "41.00"
$ws.Range("P8").Value = "41.0000"
$ws.Range("Q8").Value = "1:0"

# --- Row 9 (previously row 7): renumber to 3 ---
$ws.Range("A9").Value = 3

# --- Row 10 (previously row 8): renumber to 4 ---
$ws.Range("A10").Value = 4

# --- Row 11 (previously row 9): total selling price updates to reflect new rows ---
$ws.Range("P11").Value = 91.76

# --- Row 12 (previously row 10): update the printed timestamp ---
$ws.Range("A12").Value = "Tuesday, 5 August, 2025 9:21 AM"
